# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" sheet with the newer data pull (13 Aug 2020, 18:51).
# A handful of countries swapped rank (so the country name in a couple of rows
# changes while every row's statistics get the latest counts), and every
# numeric stat column is refreshed to the newer figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Country name swaps (rows whose ranked country changed position)
$ws.Range("A96").Value = 'Libano'
$ws.Range("A97").Value = 'Finlandia'
$ws.Range("A125").Value = 'Mozambique'
$ws.Range("A126").Value = 'Mali'
$ws.Range("A136").Value = 'Tunez'
$ws.Range("A137").Value = 'Yemen'
$ws.Range("A159").Value = 'Lesoto'
$ws.Range("A160").Value = 'Santo Tome y Principe'
$ws.Range("A161").Value = 'Aruba'
$ws.Range("A213").Value = 'Montserrat'
$ws.Range("A214").Value = 'Islas Malvinas'
$ws.Range("A1").Value = 'Datos actualizados a 13 de Agosto de 2020 a las 18:51'

# Update daily statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 5377620
$ws.Range("C4").Value = 17318
$ws.Range("D4").Value = 2815596
$ws.Range("E4").Value = 2392569
$ws.Range("G4").Value = 324
$ws.Range("H4").Value = 169455
$ws.Range("B5").Value = 3180758
$ws.Range("C5").Value = 10284
$ws.Range("E5").Value = 766753
$ws.Range("G5").Value = 265
$ws.Range("H5").Value = 104528
$ws.Range("B6").Value = 2456785
$ws.Range("C6").Value = 61314
$ws.Range("D6").Value = 1748759
$ws.Range("E6").Value = 659909
$ws.Range("G6").Value = 979
$ws.Range("H6").Value = 48117
$ws.Range("B22").Value = 221785
$ws.Range("C22").Value = 935
$ws.Range("E22").Value = 11704
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 9281
$ws.Range("B30").Value = 98343
$ws.Range("C30").Value = 1233
$ws.Range("D30").Value = 78957
$ws.Range("E30").Value = 13376
$ws.Range("G30").Value = 26
$ws.Range("H30").Value = 6010
$ws.Range("B66").Value = 28754
$ws.Range("C66").Value = 650
$ws.Range("D66").Value = 15100
$ws.Range("E66").Value = 13194
$ws.Range("G66").Value = 4
$ws.Range("H66").Value = 460
$ws.Range("B74").Value = 19286
$ws.Range("C74").Value = 211
$ws.Range("D74").Value = 13572
$ws.Range("E74").Value = 5323
$ws.Range("B96").Value = 7711
$ws.Range("C96").Value = 298
$ws.Range("D96").Value = 2496
$ws.Range("E96").Value = 5123
$ws.Range("G96").Value = 3
$ws.Range("H96").Value = 92
$ws.Range("B97").Value = 7683
$ws.Range("C97").Value = 41
$ws.Range("D97").Value = 7050
$ws.Range("E97").Value = 300
$ws.Range("H97").Value = 333
$ws.Range("B98").Value = 7368
$ws.Range("C98").Value = 68
$ws.Range("D98").Value = 6414
$ws.Range("E98").Value = 832
$ws.Range("B102").Value = 6381
$ws.Range("C102").Value = 204
$ws.Range("E102").Value = 2356
$ws.Range("G102").Value = 5
$ws.Range("H102").Value = 221
$ws.Range("D105").Value = 5167
$ws.Range("E105").Value = 132
$ws.Range("B125").Value = 2638
$ws.Range("C125").Value = 79
$ws.Range("D125").Value = 1015
$ws.Range("E125").Value = 1604
$ws.Range("H125").Value = 19
$ws.Range("B126").Value = 2582
$ws.Range("D126").Value = 1977
$ws.Range("E126").Value = 480
$ws.Range("H126").Value = 125
$ws.Range("B136").Value = 1847
$ws.Range("C136").Value = 67
$ws.Range("D136").Value = 1302
$ws.Range("E136").Value = 492
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 53
$ws.Range("B137").Value = 1841
$ws.Range("D137").Value = 937
$ws.Range("E137").Value = 376
$ws.Range("H137").Value = 528
$ws.Range("B146").Value = 1305
$ws.Range("C146").Value = 14
$ws.Range("E146").Value = 415
$ws.Range("B156").Value = 981
$ws.Range("C156").Value = 4
$ws.Range("D156").Value = 858
$ws.Range("E156").Value = 70
$ws.Range("B159").Value = 884
$ws.Range("C159").Value = 86
$ws.Range("D159").Value = 271
$ws.Range("E159").Value = 588
$ws.Range("G159").Value = 1
$ws.Range("H159").Value = 25
$ws.Range("B160").Value = 882
$ws.Range("D160").Value = 807
$ws.Range("E160").Value = 60
$ws.Range("H160").Value = 15
$ws.Range("D161").Value = 114
$ws.Range("E161").Value = 681
$ws.Range("H161").Value = 3
$ws.Range("B168").Value = 410
$ws.Range("C168").Value = 1
$ws.Range("E168").Value = 94
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
